$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# New timestamp applied to all remaining data rows
$newTimestamp = "2025-10-09 12:48:29"

# Shift rows 5 and 6 up into rows 4 and 5 (row 4's entry is removed)
$ws.Rows.Item(4).Delete()

# Update the "取得日時" timestamp for all data rows (2-5) to reflect the append time
$ws.Range("A2:A5").Value = $newTimestamp

# Column D width shrinks slightly (30 -> 28 characters)
$ws.Columns.Item(4).ColumnWidth = 28
